# Update sentiment data stats in sentiment_balance sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - German
$ws.Range("B2").Value = 1216
$ws.Range("C2").Value = 5228
$ws.Range("D2").Value = 6444
$ws.Range("E2").Value = 0.1887026691495965
$ws.Range("F2").Value = 2432

# Row 6 - Norwegian
$ws.Range("B6").Value = 5220
$ws.Range("C6").Value = 716
$ws.Range("D6").Value = 5936
$ws.Range("E6").Value = 0.8793800539083558
$ws.Range("F6").Value = 1432

# Row 11 - Chinese
$ws.Range("B11").Value = 5129
$ws.Range("C11").Value = 2797
$ws.Range("D11").Value = 7926
$ws.Range("E11").Value = 0.6471107746656574
$ws.Range("F11").Value = 5594

# Row 12 - Vietnamese
$ws.Range("B12").Value = 1880
$ws.Range("C12").Value = 2552
$ws.Range("D12").Value = 4432
$ws.Range("E12").Value = 0.4241877256317689
$ws.Range("F12").Value = 3760

# Row 14 - Cantonese
$ws.Range("B14").Value = 17897
$ws.Range("C14").Value = 18103
$ws.Range("D14").Value = 36000
$ws.Range("E14").Value = 0.4971388888888889
$ws.Range("F14").Value = 35794

# Row 16 - Finnish
$ws.Range("B16").Value = 2311
$ws.Range("C16").Value = 2175
$ws.Range("D16").Value = 4486
$ws.Range("E16").Value = 0.5151582701738743
$ws.Range("F16").Value = 4350

# Row 17 - Basque
$ws.Range("B17").Value = 1247
$ws.Range("C17").Value = 1221
$ws.Range("D17").Value = 2468
$ws.Range("E17").Value = 0.5052674230145867
$ws.Range("F17").Value = 2442

# Row 20 - Maltese
$ws.Range("B20").Value = 182
$ws.Range("C20").Value = 413
$ws.Range("D20").Value = 595
$ws.Range("E20").Value = 0.3058823529411765
$ws.Range("F20").Value = 364
